$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows at row 1017, shifting existing rows 1017-1031 down to 1025-1039
$ws.Rows.Item(1017).Resize(8).Insert()

# Row 1017
$ws.Cells.Item(1017, 1).Value = 1
$ws.Cells.Item(1017, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1017, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1017, 4).Value = 44890
$ws.Cells.Item(1017, 5).Value = 15
$ws.Cells.Item(1017, 6).Value = 100112004
$ws.Cells.Item(1017, 7).Value = "Cebolla"
$ws.Cells.Item(1017, 8).Value = "Morada(o)"
$ws.Cells.Item(1017, 9).Value = "1a (cosecha)"
$ws.Cells.Item(1017, 10).Value = 300
$ws.Cells.Item(1017, 11).Value = 5000
$ws.Cells.Item(1017, 12).Value = 5500
$ws.Cells.Item(1017, 13).Value = 5250
$ws.Cells.Item(1017, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1017, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1017, 16).Value = 292
$ws.Cells.Item(1017, 17).Value = 18
$ws.Cells.Item(1017, 18).Value = "Hortaliza"

# Row 1018
$ws.Cells.Item(1018, 1).Value = 1
$ws.Cells.Item(1018, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1018, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1018, 4).Value = 44890
$ws.Cells.Item(1018, 5).Value = 15
$ws.Cells.Item(1018, 6).Value = 100112004
$ws.Cells.Item(1018, 7).Value = "Cebolla"
$ws.Cells.Item(1018, 8).Value = "Morada(o)"
$ws.Cells.Item(1018, 9).Value = "2a (cosecha)"
$ws.Cells.Item(1018, 10).Value = 350
$ws.Cells.Item(1018, 11).Value = 4000
$ws.Cells.Item(1018, 12).Value = 4500
$ws.Cells.Item(1018, 13).Value = 4250
$ws.Cells.Item(1018, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1018, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1018, 16).Value = 236
$ws.Cells.Item(1018, 17).Value = 18
$ws.Cells.Item(1018, 18).Value = "Hortaliza"

# Row 1019
$ws.Cells.Item(1019, 1).Value = 1
$ws.Cells.Item(1019, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1019, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1019, 4).Value = 44890
$ws.Cells.Item(1019, 5).Value = 15
$ws.Cells.Item(1019, 6).Value = 100112004
$ws.Cells.Item(1019, 7).Value = "Cebolla"
$ws.Cells.Item(1019, 8).Value = "Morada(o)"
$ws.Cells.Item(1019, 9).Value = "3a (cosecha)"
$ws.Cells.Item(1019, 10).Value = 400
$ws.Cells.Item(1019, 11).Value = 3000
$ws.Cells.Item(1019, 12).Value = 3500
$ws.Cells.Item(1019, 13).Value = 3250
$ws.Cells.Item(1019, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1019, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1019, 16).Value = 181
$ws.Cells.Item(1019, 17).Value = 18
$ws.Cells.Item(1019, 18).Value = "Hortaliza"

# Row 1020
$ws.Cells.Item(1020, 1).Value = 1
$ws.Cells.Item(1020, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1020, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1020, 4).Value = 44890
$ws.Cells.Item(1020, 5).Value = 15
$ws.Cells.Item(1020, 6).Value = 100112004
$ws.Cells.Item(1020, 7).Value = "Cebolla"
$ws.Cells.Item(1020, 8).Value = "Sin especificar"
$ws.Cells.Item(1020, 9).Value = "1a (cosecha)"
$ws.Cells.Item(1020, 10).Value = 400
$ws.Cells.Item(1020, 11).Value = 9000
$ws.Cells.Item(1020, 12).Value = 10000
$ws.Cells.Item(1020, 13).Value = 9500
$ws.Cells.Item(1020, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1020, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1020, 16).Value = 528
$ws.Cells.Item(1020, 17).Value = 18
$ws.Cells.Item(1020, 18).Value = "Hortaliza"

# Row 1021
$ws.Cells.Item(1021, 1).Value = 1
$ws.Cells.Item(1021, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1021, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1021, 4).Value = 44890
$ws.Cells.Item(1021, 5).Value = 15
$ws.Cells.Item(1021, 6).Value = 100112004
$ws.Cells.Item(1021, 7).Value = "Cebolla"
$ws.Cells.Item(1021, 8).Value = "Sin especificar"
$ws.Cells.Item(1021, 9).Value = "2a (cosecha)"
$ws.Cells.Item(1021, 10).Value = 500
$ws.Cells.Item(1021, 11).Value = 7000
$ws.Cells.Item(1021, 12).Value = 8000
$ws.Cells.Item(1021, 13).Value = 7500
$ws.Cells.Item(1021, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1021, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1021, 16).Value = 417
$ws.Cells.Item(1021, 17).Value = 18
$ws.Cells.Item(1021, 18).Value = "Hortaliza"

# Row 1022
$ws.Cells.Item(1022, 1).Value = 1
$ws.Cells.Item(1022, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1022, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1022, 4).Value = 44890
$ws.Cells.Item(1022, 5).Value = 15
$ws.Cells.Item(1022, 6).Value = 100112004
$ws.Cells.Item(1022, 7).Value = "Cebolla"
$ws.Cells.Item(1022, 8).Value = "Sin especificar"
$ws.Cells.Item(1022, 9).Value = "3a (cosecha)"
$ws.Cells.Item(1022, 10).Value = 600
$ws.Cells.Item(1022, 11).Value = 6000
$ws.Cells.Item(1022, 12).Value = 7000
$ws.Cells.Item(1022, 13).Value = 6500
$ws.Cells.Item(1022, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1022, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1022, 16).Value = 361
$ws.Cells.Item(1022, 17).Value = 18
$ws.Cells.Item(1022, 18).Value = "Hortaliza"

# Row 1023
$ws.Cells.Item(1023, 1).Value = 1
$ws.Cells.Item(1023, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1023, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1023, 4).Value = 44890
$ws.Cells.Item(1023, 5).Value = 15
$ws.Cells.Item(1023, 6).Value = 100112004
$ws.Cells.Item(1023, 7).Value = "Cebolla"
$ws.Cells.Item(1023, 8).Value = "Sin especificar"
$ws.Cells.Item(1023, 9).Value = "Primera"
$ws.Cells.Item(1023, 10).Value = 400
$ws.Cells.Item(1023, 11).Value = 8000
$ws.Cells.Item(1023, 12).Value = 9000
$ws.Cells.Item(1023, 13).Value = 8500
$ws.Cells.Item(1023, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1023, 15).Value = "Perú"
$ws.Cells.Item(1023, 16).Value = 472
$ws.Cells.Item(1023, 17).Value = 18
$ws.Cells.Item(1023, 18).Value = "Hortaliza"

# Row 1024
$ws.Cells.Item(1024, 1).Value = 1
$ws.Cells.Item(1024, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1024, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(1024, 4).Value = 44890
$ws.Cells.Item(1024, 5).Value = 15
$ws.Cells.Item(1024, 6).Value = 100112004
$ws.Cells.Item(1024, 7).Value = "Cebolla"
$ws.Cells.Item(1024, 8).Value = "Sin especificar"
$ws.Cells.Item(1024, 9).Value = "Segunda"
$ws.Cells.Item(1024, 10).Value = 600
$ws.Cells.Item(1024, 11).Value = 6000
$ws.Cells.Item(1024, 12).Value = 7000
$ws.Cells.Item(1024, 13).Value = 6500
$ws.Cells.Item(1024, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(1024, 15).Value = "Perú"
$ws.Cells.Item(1024, 16).Value = 361
$ws.Cells.Item(1024, 17).Value = 18
$ws.Cells.Item(1024, 18).Value = "Hortaliza"
